# unified structure_metabolism for plants and animals
#
# This script migrates the "PLANT_BODY_DENSITY" row (row 35), which used to
# only drive plant metabolism, up into row 28 (previously
# ANIMAL_ANABOLISM_BIOMASS_CONVERSION, a row only used by animals) so that a
# single row (now named KA_ANABOLISM_FACTOR) feeds the metabolism formulas
# for both plants and animals. Row 35 is cleared out and renamed EMPTY14 to
# mark it as free for future use.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Sheet view: update frozen-pane scroll position and active selection
# ---------------------------------------------------------------------
$sheetView = $ws.Application.ActiveWindow
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollColumn = 2
$window.ScrollRow = 19
$window.RangeSelection.Worksheet.Range("H29").Select()

# ---------------------------------------------------------------------
# 2) Move the content + formatting that used to live in row 35
#    (PLANT_BODY_DENSITY) up into row 28 (A:J only - K:N belong to the
#    old ANIMAL_ANABOLISM_BIOMASS_CONVERSION row and stay as-is).
# ---------------------------------------------------------------------
$ws.Range("A35:J35").Copy($ws.Range("A28:J28"))

# Rename row 28's label and tweak the animal-specific (H) factor value.
$ws.Range("A28").Value = "KA_ANABOLISM_FACTOR"
$ws.Range("H28").Value = 0.17

# ---------------------------------------------------------------------
# 3) Clear out row 35 completely and mark it as an empty placeholder row.
#    Re-use the "red flag / empty" formatting already used by similar
#    rows (e.g. row 32) by copying its format only onto row 35.
# ---------------------------------------------------------------------
$ws.Range("A32").Copy()
$ws.Range("A35:N35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A35").Value = "EMPTY14"
$ws.Range("B35:N35").Value = 0

# ---------------------------------------------------------------------
# 4) Update the metabolism formulas (row 47, 50) that used to reference
#    row 35 so that they now reference row 28, and extend row 47's
#    ratio formula across the full K:N animal columns as well.
# ---------------------------------------------------------------------
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N")
foreach ($col in $cols) {
    $ws.Range($col + "47").Formula = "=+" + $col + "29/" + $col + "28/" + $col + "15/`$B`$3"
}

$jCols = @("B","C","D","E","F","G","H","I","J")
foreach ($col in $jCols) {
    $ws.Range($col + "50").Formula = "=+" + $col + "28*" + $col + "47*`$B`$3*" + $col + "15*" + $col + "24"
}

# Let Excel fully recalculate all dependent formulas (row 49, 50, 51, etc.)
$excel.CalculateFullRebuild()

$wb.Save()
